# Generate Report for Handback
#
# - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   across every sheet that has it.
# - The zh-cn and de-de detail sheets get two new populated columns:
#   F ("Latest Target File") and G ("Latest Handback File"), each holding a
#   hyperlinked file name, mirroring the existing A/D hyperlinked cells.
# - The "Latest Handback DateTime" (column H) timestamps are filled in with
#   real values now that the handback has happened.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    everywhere it occurs (Overview!B/C and the per-language Status column).
# ---------------------------------------------------------------------------
foreach ($sheet in $wb.Worksheets) {
    $sheet.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US", 1, [Type]::Missing, $false)
}

# ---------------------------------------------------------------------------
# Helper: find the Address of an existing hyperlink anchored at a given A1
# range (e.g. "$A$2") on a worksheet, so new hyperlinks can reuse it.
# ---------------------------------------------------------------------------
function Get-HyperlinkAddress($sheet, $targetA1) {
    foreach ($h in $sheet.Hyperlinks) {
        $rng = $h.Range()
        $a1 = $rng.Address()
        if ($a1 -eq $targetA1) {
            return $h.Address()
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 2) Populate "Latest Target File" (F) / "Latest Handback File" (G) for rows
#    2 and 3 on both the zh-cn and de-de sheets, and fill in the real
#    "Latest Handback DateTime" (H) now that handback has occurred.
# ---------------------------------------------------------------------------
$langSheets = @(
    @{ Name = "zh-cn"; HandbackDate = "2016-03-20 02:25:37" },
    @{ Name = "de-de"; HandbackDate = "2016-03-20 02:25:43" }
)

foreach ($entry in $langSheets) {
    $ws = $wb.Worksheets.Item($entry.Name)

    foreach ($row in 2, 3) {
        $aAddr = '$A$' + $row
        $dAddr = '$D$' + $row

        $mdAddress = Get-HyperlinkAddress $ws $aAddr
        $xlfAddress = Get-HyperlinkAddress $ws $dAddr

        $aCell = $ws.Range("A$row")
        $dCell = $ws.Range("D$row")
        $fCell = $ws.Range("F$row")
        $gCell = $ws.Range("G$row")

        $mdDisplay = $aCell.Value()
        $xlfDisplay = $dCell.Value()

        $fCell.Value = $mdDisplay
        if ($mdAddress) {
            $ws.Hyperlinks.Add($fCell, $mdAddress, [Type]::Missing, [Type]::Missing, $mdDisplay) | Out-Null
        }

        $gCell.Value = $xlfDisplay
        if ($xlfAddress) {
            $ws.Hyperlinks.Add($gCell, $xlfAddress, [Type]::Missing, [Type]::Missing, $xlfDisplay) | Out-Null
        }

        # Latest Handback DateTime
        $ws.Range("H$row").Value = $entry.HandbackDate
    }
}

"Handback report generated"
